$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -9.006899999999996
$ws.Range("D9").Value = -7.775799999999993
$ws.Range("A11").Value = -21.5366
$ws.Range("B11").Value = 5.833900000000001
$ws.Range("A12").Value = -21.3993
$ws.Range("D13").Value = -8.808199999999996
$ws.Range("D14").Value = -8.3935
$ws.Range("A15").Value = -21.7258
$ws.Range("D19").Value = -7.933200000000002
$ws.Range("D21").Value = -8.161299999999994
$ws.Range("D22").Value = -7.983199999999999
$ws.Range("B23").Value = 8.541899999999996
$ws.Range("D24").Value = -7.980599999999996
$ws.Range("D26").Value = -8.397399999999999
$ws.Range("A27").Value = -21.9426
$ws.Range("A28").Value = -21.8268
$ws.Range("B28").Value = 6.066400000000002
$ws.Range("A31").Value = -21.77960000000001
$ws.Range("A32").Value = -21.04829999999999
$ws.Range("B32").Value = 5.4677
$ws.Range("B34").Value = 9.598700000000003
$ws.Range("A36").Value = -21.13209999999999
$ws.Range("B36").Value = 5.064899999999998
$ws.Range("B37").Value = 8.945200000000005
$ws.Range("A38").Value = -20.1129
$ws.Range("D38").Value = -7.732099999999996
$ws.Range("D41").Value = -8.569399999999996
$ws.Range("B42").Value = 10.232
$ws.Range("A46").Value = -22.01710000000001
$ws.Range("B49").Value = 4.662899999999999
$ws.Range("D52").Value = -8.179800000000002
$ws.Range("A54").Value = -21.8087
$ws.Range("B54").Value = 5.832599999999996
$ws.Range("A55").Value = -21.76760000000002
$ws.Range("A56").Value = -21.82139999999999
$ws.Range("D56").Value = -8.948199999999998
$ws.Range("A67").Value = -21.55519999999996
$ws.Range("A69").Value = -21.72309999999998
$ws.Range("D71").Value = -7.307299999999997
$ws.Range("A72").Value = -21.9134
$ws.Range("D72").Value = -7.426500000000005
$ws.Range("A73").Value = -19.262
$ws.Range("B78").Value = 8.615399999999996
$ws.Range("D78").Value = -8.394100000000002
$ws.Range("B80").Value = 9.734999999999998
$ws.Range("A83").Value = -21.88529999999999
$ws.Range("D83").Value = -8.4735
$ws.Range("D85").Value = -9.066499999999996
$ws.Range("A86").Value = -21.98110000000002
$ws.Range("D86").Value = -8.5221
$ws.Range("D90").Value = -6.776499999999996
$ws.Range("A91").Value = -20.71499999999999
$ws.Range("A93").Value = -21.54830000000002
$ws.Range("D96").Value = -8.821899999999996
$ws.Range("B97").Value = 6.605399999999999
$ws.Range("A99").Value = -21.7917
$ws.Range("B99").Value = 6.488499999999998
$ws.Range("B100").Value = 4.3519
$ws.Range("B101").Value = 4.9954
$ws.Range("D103").Value = -8.431499999999998
$ws.Range("A104").Value = -21.4462
$ws.Range("A105").Value = -19.92609999999998
